# Update cryptos list Price (D) and Volume(1h) (E) columns
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text storage so numeric-looking strings (e.g. "0.9995") are not
# reinterpreted as numbers, matching the original inline-string cell type.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "30.317.47"
$ws.Range("E2").Value = "  -1.17%  "
$ws.Range("D3").Value = "1.868.99"
$ws.Range("E3").Value = "  -0.74%  "
$ws.Range("D4").Value = "0.9995"
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "243.85"
$ws.Range("E5").Value = "  -2.31%  "
$ws.Range("D6").Value = "0.9997"
$ws.Range("E6").Value = "  -0.02%  "
$ws.Range("D7").Value = "0.4726"
$ws.Range("E7").Value = "  -0.44%  "
$ws.Range("D8").Value = "0.2869"
$ws.Range("E8").Value = "  -2.15%  "
$ws.Range("D9").Value = "0.06482"
$ws.Range("E9").Value = "  -0.84%  "
$ws.Range("D10").Value = "21.74"
$ws.Range("E10").Value = "  -1.31%  "
$ws.Range("D11").Value = "0.07800"
$ws.Range("E11").Value = "  +0.96%  "
$ws.Range("D12").Value = "98.61"
$ws.Range("E12").Value = "  +1.93%  "
$ws.Range("D13").Value = "1.867.72"
$ws.Range("E13").Value = "  -0.85%  "
$ws.Range("D14").Value = "0.7244"
$ws.Range("E14").Value = "  -2.08%  "
$ws.Range("D15").Value = "5.161"
$ws.Range("E15").Value = "  -1.52%  "
$ws.Range("D16").Value = "281.40"
$ws.Range("E16").Value = "  +2.09%  "
$ws.Range("D17").Value = "30.297.06"
$ws.Range("E17").Value = "  -1.19%  "
$ws.Range("D18").Value = "13.06"
$ws.Range("E18").Value = "  -0.93%  "
$ws.Range("D19").Value = "0.9997"
$ws.Range("D20").Value = "0.000007470"
$ws.Range("E20").Value = "  -1.19%  "
$ws.Range("D21").Value = "2.111.98"
$ws.Range("E21").Value = "  -0.90%  "
$ws.Range("D22").Value = "0.9998"
$ws.Range("E22").Value = "  -0.04%  "
$ws.Range("D23").Value = "5.253"
$ws.Range("E23").Value = "  -1.49%  "
$ws.Range("D24").Value = "6.273"
$ws.Range("E24").Value = "  +0.66%  "
$ws.Range("D25").Value = "162.48"
$ws.Range("E25").Value = "  -1.04%  "
$ws.Range("D26").Value = "9.032"
$ws.Range("E26").Value = "  -2.23%  "
$ws.Range("D27").Value = "18.76"
$ws.Range("E27").Value = "  -0.37%  "
$ws.Range("D28").Value = "1.884"
$ws.Range("E28").Value = "  -1.81%  "
$ws.Range("D29").Value = "0.09641"
$ws.Range("E29").Value = "  -0.70%  "
$ws.Range("D30").Value = "1.314"
$ws.Range("E30").Value = "  -2.30%  "
$ws.Range("D31").Value = "1.481"
$ws.Range("E31").Value = "  -1.63%  "
$ws.Range("D32").Value = "4.229"
$ws.Range("E32").Value = "  -1.50%  "
$ws.Range("D33").Value = "4.128"
$ws.Range("E33").Value = "  -0.89%  "
$ws.Range("D34").Value = "0.04797"
$ws.Range("E34").Value = "  -2.22%  "
$ws.Range("D35").Value = "1.120"
$ws.Range("E35").Value = "  -0.51%  "
$ws.Range("D36").Value = "0.6869"
$ws.Range("E36").Value = "  -1.92%  "
$ws.Range("D37").Value = "2.713"
$ws.Range("E37").Value = "  -0.54%  "
$ws.Range("D38").Value = "0.01891"
$ws.Range("E38").Value = "  -0.97%  "
$ws.Range("D39").Value = "2.830"
$ws.Range("E39").Value = "  +1.36%  "
$ws.Range("D40").Value = "75.34"
$ws.Range("E40").Value = "  -0.53%  "
$ws.Range("D41").Value = "6.241"
$ws.Range("E41").Value = "  -1.39%  "
$ws.Range("D42").Value = "1.944"
$ws.Range("E42").Value = "  -3.85%  "
$ws.Range("D43").Value = "0.4216"
$ws.Range("E43").Value = "  -0.75%  "
$ws.Range("D44").Value = "0.9988"
$ws.Range("E44").Value = "  -0.10%  "
$ws.Range("D45").Value = "0.8232"
$ws.Range("E45").Value = "  -2.07%  "
$ws.Range("D46").Value = "100.71"
$ws.Range("E46").Value = "  -1.88%  "
$ws.Range("D47").Value = "9.635"
$ws.Range("E47").Value = "  +2.29%  "
$ws.Range("D48").Value = "6.987"
$ws.Range("E48").Value = "  -0.94%  "
$ws.Range("D49").Value = "35.03"
$ws.Range("E49").Value = "  -1.72%  "
$ws.Range("D50").Value = "0.05765"
$ws.Range("D51").Value = "882.86"
$ws.Range("E51").Value = "  -3.96%  "

# Restore default (unstyled) cell style so no stray number-format survives
# on the cells, matching the original workbook formatting.
$ws.Range("D2:E51").Style = "Normal"
